# feat: import jp per pekan excel
#
# Adds a new "Jatah Per Pekan" column to the Table1 listobject on Sheet1,
# between the existing "Jenis Mapel" column (B) and the free-standing
# "Mata Pelajaran" header (previously C, now shifted to D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank column at C - this shifts the existing "Mata Pelajaran"
# header (and its column formatting) from C to D.
$ws.Columns.Item(3).Insert()

# Grow Table1 (currently A1:B14) one column to the right so it covers the
# freshly inserted column C.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:C14"))

# Name the new table column / header cell (must happen after Resize, which
# otherwise auto-names the new column "Column3").
$ws.Range("C1").Value = "Jatah Per Pekan"

# Match the column widths the author ended up with after inserting the
# column and nudging the borders.
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 27

# Leave the selection where the author left it, on the (now relocated)
# "Mata Pelajaran" column.
$ws.Range("D5").Select()
